$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsShare = $wb.Worksheets.Item("SoDSCbRIC")

# -----------------------------------------------------------------
# 1) Data sheet: split the "ISIC 20T21 / Chemicals and pharmaceutical
#    products" lookup row (G11:H11) into two rows:
#       G11:H11 -> "ISIC 20" / "Chemicals"
#       G12:H12 -> "ISIC 21" / "Pharmaceuticals"
#    Only columns G:H of the lookup table move - columns A:E (the
#    separate cost-breakdown table on the same sheet) must stay put.
#    Shift existing G12:H37 down to G13:H38 first (bottom-up so we
#    don't clobber values before they are copied).
# -----------------------------------------------------------------
for ($r = 37; $r -ge 12; $r--) {
    $gVal = $wsData.Cells.Item($r, 7).Value()
    $hVal = $wsData.Cells.Item($r, 8).Value()
    $wsData.Cells.Item($r + 1, 7).Value = $gVal
    $wsData.Cells.Item($r + 1, 8).Value = $hVal
}

$wsData.Range("G11").Value = "ISIC 20"
$wsData.Range("H11").Value = "Chemicals"
$wsData.Range("G12").Value = "ISIC 21"
$wsData.Range("H12").Value = "Pharmaceuticals"

# -----------------------------------------------------------------
# 2) SoDSCbRIC sheet: insert a new column at L (shifting L:AK right
#    to M:AL) to hold the new "ISIC 21" bucket, keep K1 as "ISIC 20".
# -----------------------------------------------------------------
$wsShare.Range("L1").EntireColumn.Insert(-4161)  # xlShiftToRight

$wsShare.Range("K1").Value = "ISIC 20"
$wsShare.Range("L1").Value = "ISIC 21"

$wsShare.Range("L2").Formula = "=SUMIF(Data!`$E`$3:`$E`$13,SoDSCbRIC!L1,Data!`$D`$3:`$D`$13)"
$wsShare.Range("L3").Formula = "=L2"
$wsShare.Range("L4").Formula = "=SUMIF(Data!`$E`$17:`$E`$27,SoDSCbRIC!L1,Data!`$D`$17:`$D`$27)"
